$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the existing row 860 (2026/02/24, 2:00 entry).
# This shifts the old rows 861..902 down to 863..904, preserving their values.
$ws.Rows.Item(861).EntireRow.Insert()
$ws.Rows.Item(861).EntireRow.Insert()

# Fill the two newly inserted rows with the additional 2026/02/24 time slots.
# The date strings are stored as plain text (like the rest of the sheet), so a
# leading apostrophe stops Excel from auto-converting them into date serials;
# resetting the style back to Normal afterwards removes the quote-prefix style
# so the cell matches the plain (unstyled) cells used throughout the sheet.

$ws.Cells.Item(861, 1).Value = "'2026/02/24"
$ws.Cells.Item(861, 1).Style = "Normal"
$ws.Cells.Item(861, 2).Value = "火"
$ws.Cells.Item(861, 3).Value = 7
$ws.Cells.Item(861, 4).Value = 201

$ws.Cells.Item(862, 1).Value = "'2026/02/24"
$ws.Cells.Item(862, 1).Style = "Normal"
$ws.Cells.Item(862, 2).Value = "火"
$ws.Cells.Item(862, 3).Value = 10
$ws.Cells.Item(862, 4).Value = 201
